$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header of column G: it was incorrectly labeled "datetime_eri_b"
# (duplicate of column A) but should read "datetime_eri_c" since columns
# G/H hold the datetime/band for star C.
$ws.Range("G1").Value = "datetime_eri_c"

# Update the active cell selection to match the saved state.
$ws.Range("F15").Select()
